# Weekly update: insert this week's "Betarraga" (Femacal de La Calera) price
# records at the top of the data block (rows 422-423), pushing the existing
# history down by two rows. The oldest week (rows 420-421) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 422-423; everything below (old 422:438) shifts
# down to 424:440, carrying its formatting (incl. the date style on column D).
$ws.Range("A422:A423").EntireRow.Insert()

# --- New row 422: Betarraga, "Primera" quality, week of 2021-11-09 ---
$ws.Range("A422").Value = 3
$ws.Range("B422").Value = "Femacal de La Calera"
$ws.Range("C422").Value = "Coquimbo"
$ws.Range("D422").Value = 44509
$ws.Range("E422").Value = 5
$ws.Range("F422").Value = 100114014
$ws.Range("G422").Value = "Betarraga"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 3600
$ws.Range("K422").Value = 500
$ws.Range("L422").Value = 550
$ws.Range("M422").Value = 525
$ws.Range("N422").Value = "$/paquete 4 unidades"
$ws.Range("O422").Value = "Provincia de Quillota"
$ws.Range("P422").Value = 131
$ws.Range("Q422").Value = 4
$ws.Range("R422").Value = "Hortaliza"

# --- New row 423: Betarraga, "Segunda" quality, week of 2021-11-09 ---
$ws.Range("A423").Value = 3
$ws.Range("B423").Value = "Femacal de La Calera"
$ws.Range("C423").Value = "Coquimbo"
$ws.Range("D423").Value = 44509
$ws.Range("E423").Value = 5
$ws.Range("F423").Value = 100114014
$ws.Range("G423").Value = "Betarraga"
$ws.Range("H423").Value = "Sin especificar"
$ws.Range("I423").Value = "Segunda"
$ws.Range("J423").Value = 1600
$ws.Range("K423").Value = 400
$ws.Range("L423").Value = 400
$ws.Range("M423").Value = 400
$ws.Range("N423").Value = "$/paquete 4 unidades"
$ws.Range("O423").Value = "Provincia de Quillota"
$ws.Range("P423").Value = 100
$ws.Range("Q423").Value = 4
$ws.Range("R423").Value = "Hortaliza"
